$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '58.798.68'
$ws.Range('E2').Value = '  -2.95%  '
$ws.Range('D3').Value = '2.724.65'
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '505.13'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -4.48%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '141.17'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -1.23%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.997'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.16%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.531'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -4.10%  '
$ws.Range('D9').Value = '2.737.05'
$ws.Range('E9').Value = '  -5.99%  '
$ws.Range('E10').Value = '  -2.48%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '6.07'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +2.05%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.349'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -3.02%  '
$ws.Range('E13').Value = '  +0.95%  '
$ws.Range('D14').Value = '3.199.46'
$ws.Range('D15').Value = '58.813.23'
$ws.Range('E15').Value = '  -2.99%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '21.66'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -4.33%  '
$ws.Range('B17').Value = 'WrappedEther'
$ws.Range('C17').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D17').Value = '2.735.86'
$ws.Range('E17').Value = '  -5.84%  '
$ws.Range('B18').Value = 'ShibaInu'
$ws.Range('C18').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.0000136'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -4.47%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '4.77'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -5.47%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '10.98'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -6.18%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '341.83'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -6.04%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.26'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -5.77%  '
$ws.Range('E23').Value = '  -0.24%  '
$ws.Range('E24').Value = '  -0.04%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '63.38'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -1.22%  '
$ws.Range('E26').Value = '  -2.50%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.427'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -5.67%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.994'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.52%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.51'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -3.85%  '
$ws.Range('D30').Value = '0.0₃0829'
$ws.Range('E30').Value = '  -3.91%  '
$ws.Range('E31').Value = '  -0.04%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '19.24'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -2.09%  '
$ws.Range('E33').Value = '  -4.19%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '150.87'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +2.03%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '4.21'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -3.61%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '5.43'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -2.89%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.947'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -5.85%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.13'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -6.06%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '36.29'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -4.30%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.59'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -2.27%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.39'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -7.15%  '
$ws.Range('D42').Value = '2.193.46'
$ws.Range('E42').Value = '  -5.97%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.0563'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -2.23%  '
$ws.Range('E44').Value = '  -0.13%  '
$ws.Range('E45').Value = '  -6.72%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '19.05'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -8.05%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '4.78'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -3.71%  '
$ws.Range('E48').Value = '  +0.25%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0227'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -3.14%  '
$ws.Range('E50').Value = '  -4.60%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '18.08'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -2.09%  '
